$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175898432731628
$ws.Range("B1").Value = 1.416824102401733
$ws.Range("C1").Value = 1.250287532806396
$ws.Range("D1").Value = 1.362024784088135
$ws.Range("E1").Value = 1.22085440158844
